$d = $word.ActiveDocument

$replacements = @(
    @("152×4=", "553×5="),
    @("288×9=", "523×3="),
    @("220×2=", "483×3="),
    @("249×3=", "720×3="),
    @("387×2=", "208×6="),
    @("470×6=", "555×3="),
    @("851×7=", "427×4="),
    @("702×5=", "843×6="),
    @("654×5=", "633×3="),
    @("647×5=", "369×9="),
    @("379×3=", "649×2="),
    @("106×7=", "551×8="),
    @("641×3=", "938×4="),
    @("120×6=", "156×4="),
    @("370×7=", "731×4="),
    @("448×8=", "624×3="),
    @("541×7=", "167×5="),
    @("532×7=", "529×5="),
    @("393×3=", "575×2="),
    @("947×8=", "165×8="),
    @("406×3=", "736×7="),
    @("157×2=", "703×6="),
    @("704×6=", "440×5="),
    @("581×9=", "191×6="),
    @("926×6=", "528×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
